$d = $word.ActiveDocument

# boundary 0
$ok0 = $d.Content.Find.Execute(" McGraw-Hill/Irwin. 2021.Crawford, M., Di Benedett", $true, $false, $false, $false, $false, $true, 1, $false, " McGraw-Hill/Irwin. 2021.^lCrawford, M., Di Benedett", 2)
Write-Host "boundary 0: $ok0"

# boundary 1
$ok1 = $d.Content.Find.Execute(" McGraw Hill Brasil, 2016uma referência para a mel", $true, $false, $false, $false, $false, $true, 1, $false, " McGraw Hill Brasil, 2016^luma referência para a mel", 2)
Write-Host "boundary 1: $ok1"

# boundary 2
$ok2 = $d.Content.Find.Execute("lo: Saraiva, 2006. 542p. CHENG, L.C. QFD: desdobra", $true, $false, $false, $false, $false, $true, 1, $false, "lo: Saraiva, 2006. 542p. ^lCHENG, L.C. QFD: desdobra", 2)
Write-Host "boundary 2: $ok2"

# boundary 3
$ok3 = $d.Content.Find.Execute("o: Edgard Blücher, 2007. DENIS, R. C. Uma introduç", $true, $false, $false, $false, $false, $true, 1, $false, "o: Edgard Blücher, 2007. ^lDENIS, R. C. Uma introduç", 2)
Write-Host "boundary 3: $ok3"

# boundary 4
$ok4 = $d.Content.Find.Execute(": Edgard Blücher., 2000. MALHOTRA, N.K. Pesquisa d", $true, $false, $false, $false, $false, $true, 1, $false, ": Edgard Blücher., 2000. ^lMALHOTRA, N.K. Pesquisa d", 2)
Write-Host "boundary 4: $ok4"

# boundary 5
$ok5 = $d.Content.Find.Execute("o Alegre: Bookman, 2006. PAHL,G.; BEITZ,W.; FELDHU", $true, $false, $false, $false, $false, $true, 1, $false, "o Alegre: Bookman, 2006. ^lPAHL,G.; BEITZ,W.; FELDHU", 2)
Write-Host "boundary 5: $ok5"

# boundary 6
$ok6 = $d.Content.Find.Execute("lo: Edgard Blucher, 2005.ROZENFELD, H.; FORCELLINI", $true, $false, $false, $false, $false, $true, 1, $false, "lo: Edgard Blucher, 2005.^lROZENFELD, H.; FORCELLINI", 2)
Write-Host "boundary 6: $ok6"

